# Atualiza instrucao de trabalho
# - Marks several existing incident rows as "Resolvido" (previously "Pendente")
# - Appends new pending incident rows to both the SPN and ITI sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: SPN
# ---------------------------------------------------------------------------
$wsSPN = $wb.Worksheets.Item("SPN")

# Update status column (J) from "Pendente" to "Resolvido" for rows 149-153
$spnResolvedRows = @(149, 150, 151, 152, 153)
foreach ($r in $spnResolvedRows) {
    $wsSPN.Cells.Item($r, 10).Value = "Resolvido"
}

# Append new incident rows 154-157
$spnNewRows = @(
    @(154, "SPN", "Higor Cruz",  2025, 33, "18/08/2025", "22/08/2025", 343563, "08/2025", "18/08/2025", "Pendente", "Willian Rios"),
    @(155, "SPN", "Higor Cruz",  2025, 33, "18/08/2025", "22/08/2025", 343881, "08/2025", "18/08/2025", "Pendente", "Willian Rios"),
    @(156, "SPN", "Luan Pierry", 2025, 33, "18/08/2025", "22/08/2025", 343975, "08/2025", "18/08/2025", "Pendente", "Willian Rios"),
    @(157, "SPN", "Mara Neves",  2025, 33, "18/08/2025", "22/08/2025", 343106, "08/2025", "18/08/2025", "Pendente", "Willian Rios")
)

foreach ($row in $spnNewRows) {
    $r = $row[0]
    $wsSPN.Cells.Item($r, 1).Value  = $row[1]
    $wsSPN.Cells.Item($r, 2).Value  = $row[2]
    $wsSPN.Cells.Item($r, 3).Value  = $row[3]
    $wsSPN.Cells.Item($r, 4).Value  = $row[4]
    $wsSPN.Cells.Item($r, 5).Value  = $row[5]
    $wsSPN.Cells.Item($r, 6).Value  = $row[6]
    $wsSPN.Cells.Item($r, 7).Value  = $row[7]
    $wsSPN.Cells.Item($r, 8).Value  = $row[8]
    $wsSPN.Cells.Item($r, 9).Value  = $row[9]
    $wsSPN.Cells.Item($r, 10).Value = $row[10]
    $wsSPN.Cells.Item($r, 11).Value = $row[11]
}

# ---------------------------------------------------------------------------
# Sheet: ITI
# ---------------------------------------------------------------------------
$wsITI = $wb.Worksheets.Item("ITI")

# Update status column (J) from "Pendente" to "Resolvido" for the listed rows
$itiResolvedRows = @(359, 392, 410, 413, 418, 419, 420, 421, 422, 425, 426, 427, 428, 429)
foreach ($r in $itiResolvedRows) {
    $wsITI.Cells.Item($r, 10).Value = "Resolvido"
}

# Append new incident rows 433-442
$itiNewRows = @(
    @(433, "ITI", "Erick da Silva",    2025, 33, "18/08/2025", "22/08/2025", 343419, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(434, "ITI", "Erick da Silva",    2025, 33, "18/08/2025", "22/08/2025", 343497, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(435, "ITI", "Gabriel López",     2025, 33, "18/08/2025", "22/08/2025", 343928, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(436, "ITI", "Guilherme Worel",   2025, 33, "18/08/2025", "22/08/2025", 344118, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(437, "ITI", "Guilherme Worel",   2025, 33, "18/08/2025", "22/08/2025", 343917, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(438, "ITI", "Guilherme Worel",   2025, 33, "18/08/2025", "22/08/2025", 343649, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(439, "ITI", "Lourival Moizés",   2025, 33, "18/08/2025", "22/08/2025", 343933, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(440, "ITI", "Sostenes Simões",   2025, 33, "18/08/2025", "22/08/2025", 343986, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(441, "ITI", "Sostenes Simões",   2025, 33, "18/08/2025", "22/08/2025", 344125, "08/2025", "18/08/2025", "Pendente", "Emerson Simette"),
    @(442, "ITI", "Sostenes Simões",   2025, 33, "18/08/2025", "22/08/2025", 344168, "08/2025", "18/08/2025", "Pendente", "Emerson Simette")
)

foreach ($row in $itiNewRows) {
    $r = $row[0]
    $wsITI.Cells.Item($r, 1).Value  = $row[1]
    $wsITI.Cells.Item($r, 2).Value  = $row[2]
    $wsITI.Cells.Item($r, 3).Value  = $row[3]
    $wsITI.Cells.Item($r, 4).Value  = $row[4]
    $wsITI.Cells.Item($r, 5).Value  = $row[5]
    $wsITI.Cells.Item($r, 6).Value  = $row[6]
    $wsITI.Cells.Item($r, 7).Value  = $row[7]
    $wsITI.Cells.Item($r, 8).Value  = $row[8]
    $wsITI.Cells.Item($r, 9).Value  = $row[9]
    $wsITI.Cells.Item($r, 10).Value = $row[10]
    $wsITI.Cells.Item($r, 11).Value = $row[11]
}
